$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.911.93'
$ws.Range("E2").Value = '  -0.73%  '

# Row 3
$ws.Range("D3").Value = '1.741.43'
$ws.Range("E3").Value = '  +1.49%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.07%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5070'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.12%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3577'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.38%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.10'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.08%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07233'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.40%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.061'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.09%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.05%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.15'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.15%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.949'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.75%  '

# Row 15
$ws.Range("D15").Value = '1.739.90'
$ws.Range("E15").Value = '  +1.65%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.801'
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.20%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001031'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.82%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06404'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.76%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9997'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.06%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.06%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.747'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.03%  '

# Row 23
$ws.Range("D23").Value = '26.979.91'
$ws.Range("E23").Value = '  -0.68%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.35%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.037'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.44%  '

# Row 26
$ws.Range("E26").Value = '  -3.13%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.75%  '

# Row 28
$ws.Range("D28").Value = '1.939.09'
$ws.Range("E28").Value = '  +1.46%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.217'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.24%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.20%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.038'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.95%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09597'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.79%  '

# Row 33
$ws.Range("E33").Value = '  -0.25%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.371'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.03%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05895'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.43%  '

# Row 36
$ws.Range("E36").Value = '  -0.65%  '

# Row 37
$ws.Range("E37").Value = '  +0.06%  '

# Row 38
$ws.Range("E38").Value = '  +0.35%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.423'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.20%  '

# Row 40
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.756'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.46%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6030'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.103'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.84%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.621'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.46%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.86%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.586'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.96%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5647'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.03%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '120.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.87%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.844'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.44%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.099'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.23%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06647'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.01%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9997'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.04%  '
